$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Date column (B) with new timestamps for rows 2-7
$ws.Range("B2").Value = "Tue Jan 28 22:09:42 EST 2025"
$ws.Range("B3").Value = "Tue Jan 28 22:09:55 EST 2025"
$ws.Range("B4").Value = "Tue Jan 28 22:10:07 EST 2025"
$ws.Range("B5").Value = "Tue Jan 28 22:10:19 EST 2025"
$ws.Range("B6").Value = "Tue Jan 28 22:10:31 EST 2025"
$ws.Range("B7").Value = "Tue Jan 28 22:10:43 EST 2025"

# Update rows 6-7: Result changes from Fail to Pass, Execute changes from DoNotRun to Y
$ws.Range("A6").Value = "Pass"
$ws.Range("C6").Value = "Y"

$ws.Range("A7").Value = "Pass"
$ws.Range("C7").Value = "Y"

# Update selection to C7 only (was C6:C7 with active cell C6)
$null = $ws.Range("C7").Select()
